# The scraper re-ran and appended/refreshed its snapshot at 2025-09-09 18:34 JST,
# so every existing row's "取得日時" (fetched-at) timestamp in column A is bumped
# from 18:21:53 to 18:34:52 on the "ランサーズ" (案件情報) sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-09 18:34:52"

$ws.Range("A2:A15").Value = $newTimestamp
